# "Ya no nos falta tanto"
# - Remove the finished "CATALOGO..." task row (old row 24), shifting the rows below it up.
# - Mark three tasks ("Usuario dividido...", "Implementar la sesion...", "Arreglar botones...")
#   as assigned to Araujo, highlighted in dark red.
# - Mark several other tasks as already done ("x") with the red-highlighted style centered.
# - Mark the ""otro index"" task as no longer needed ("Este ya no se hace"), highlighted in red.
# - Center-align the Hecho/No hecho columns (D:E) across the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Delete the completed "CATALOGO" task row entirely (old row 24).
#    Everything below shifts up by one row.
# ------------------------------------------------------------------
$ws.Rows("24:24").Delete()

# ------------------------------------------------------------------
# 2. New "assigned to Araujo" highlight block (rows 6-8), dark red fill + border.
# ------------------------------------------------------------------
$ws.Range("C6").Value = "Araujo"
$ws.Range("D6").Value = "x"

$ws.Range("C7").Value = "Araujo"
$ws.Range("D7").Value = "x"

# D8 already contains "x"; C8 stays empty.

$ws.Range("B6:C8").Interior.Color = 192
$ws.Range("B6:C8").Borders.LineStyle = 1

$ws.Range("D6:E8").Interior.Color = 192
$ws.Range("D6:E8").Borders.LineStyle = 1
$ws.Range("D6:E8").HorizontalAlignment = -4108

# Stray centered (no fill/border) cell that shows up next to row 7 in the target file.
$ws.Range("F7").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 3. Mark additional "x" (hecho) cells that were blank before.
# ------------------------------------------------------------------
$ws.Range("D9").Value = "x"
$ws.Range("D10").Value = "x"
$ws.Range("D12").Value = "x"
$ws.Range("D15").Value = "x"

# ------------------------------------------------------------------
# 4. Row 22 ("otro index" task): no longer needed, highlighted in red.
# ------------------------------------------------------------------
$ws.Range("C22").Value = "Este ya no se hace"
$ws.Range("D22").Value = "x"
$ws.Range("E22").Value = "x"

$ws.Range("B22:C22").Interior.Color = 255
$ws.Range("B22:C22").Borders.LineStyle = 1

$ws.Range("D22:E22").Interior.Color = 255
$ws.Range("D22:E22").Borders.LineStyle = 1
$ws.Range("D22:E22").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 5. Center the Hecho (D) / No hecho (E) columns across the whole table
#    (rows that still use the plain bordered style get a border + center;
#    rows already filled red/white just pick up the centering).
# ------------------------------------------------------------------
$ws.Range("D5:E5").HorizontalAlignment = -4108
$ws.Range("D9:E9").HorizontalAlignment = -4108
$ws.Range("D10:E10").HorizontalAlignment = -4108
$ws.Range("D11:E11").HorizontalAlignment = -4108
$ws.Range("D12:E12").HorizontalAlignment = -4108
$ws.Range("D13:E13").HorizontalAlignment = -4108
$ws.Range("D14:E14").HorizontalAlignment = -4108
$ws.Range("D15:E15").HorizontalAlignment = -4108
$ws.Range("D16:E16").HorizontalAlignment = -4108
$ws.Range("D17:E17").HorizontalAlignment = -4108
$ws.Range("D18:E18").HorizontalAlignment = -4108
$ws.Range("D19:E19").HorizontalAlignment = -4108
$ws.Range("D20:E20").HorizontalAlignment = -4108
$ws.Range("D21:E21").HorizontalAlignment = -4108
$ws.Range("D23:E23").HorizontalAlignment = -4108
$ws.Range("D24:E24").HorizontalAlignment = -4108
$ws.Range("E25").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 6. View / window cosmetics.
# ------------------------------------------------------------------
$ws.Range("E5:E25").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
